# Update "想去人数" (number of people interested) figures across the
# workbook's four sheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4679
$ws1.Range("F5").Value  = 183
$ws1.Range("F6").Value  = 1853
$ws1.Range("F8").Value  = 753
$ws1.Range("F9").Value  = 38
$ws1.Range("F11").Value = 417
$ws1.Range("F12").Value = 1136
$ws1.Range("F13").Value = 1589
$ws1.Range("F15").Value = 1658
$ws1.Range("F16").Value = 559
$ws1.Range("F17").Value = 523
$ws1.Range("F19").Value = 189
$ws1.Range("F20").Value = 1554
$ws1.Range("F21").Value = 1194
$ws1.Range("F22").Value = 599
$ws1.Range("F23").Value = 2517
$ws1.Range("F26").Value = 1584
$ws1.Range("F30").Value = 71
$ws1.Range("F31").Value = 4278

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 4164
$ws2.Range("F17").Value = 283
$ws2.Range("F29").Value = 3
$ws2.Range("F38").Value = 37

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1341
$ws3.Range("F5").Value = 1735
$ws3.Range("F7").Value = 292

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1341
$ws4.Range("F3").Value  = 1735
$ws4.Range("F5").Value  = 292
$ws4.Range("F9").Value  = 4679
$ws4.Range("F11").Value = 183
$ws4.Range("F12").Value = 1853
$ws4.Range("F14").Value = 753
$ws4.Range("F15").Value = 38
$ws4.Range("F19").Value = 417
$ws4.Range("F20").Value = 1136
$ws4.Range("F21").Value = 1589
$ws4.Range("F24").Value = 1658
$ws4.Range("F25").Value = 559
$ws4.Range("F26").Value = 523
$ws4.Range("F28").Value = 189
$ws4.Range("F29").Value = 283
$ws4.Range("F32").Value = 1554
$ws4.Range("F33").Value = 1194
$ws4.Range("F34").Value = 599
$ws4.Range("F37").Value = 2517
$ws4.Range("F43").Value = 1584
$ws4.Range("F47").Value = 71
$ws4.Range("F48").Value = 4278
$ws4.Range("F49").Value = 37
